$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New kCar rows to add to PlotShearFlow fit results sheet
$data = @(
    @("10%_0WSt_kCar_1", -164.5681764470073, 3.645852121809451, 0.07752968282954108, 0.002211246080797293, 404.0117296072254, 3.491151435698671),
    @("10%_0WSt_kCar_2", -130.5920206515675, 3.223792544422935, 0.08099504891644603, 0.002479759855658115, 390.4951195188706, 3.098563056362585),
    @("10%_0WSt_kCar_3", -62.0559463280163, 2.988344031710459, 0.1922832153837944, 0.006000398964713743, 541.0181235241719, 3.373781092379696)
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
